$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text updates (rich-text shared strings): volume/issue number
# and the reporting week date range.
# ------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "18"
$ws.Range("C9").Characters(27, 9).Text = "4/28/2025"
$ws.Range("C9").Characters(47, 9).Text = "5/4/2025"

# ------------------------------------------------------------------
# Precinct crime-complaint statistics table (rows 14-31): refreshed
# weekly/28-day/YTD/2-year counts and recomputed % changes.
# ------------------------------------------------------------------
# Row 14
$ws.Range("N14").Value = -75

# Row 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 10
$ws.Range("L15").Value = 57.142857142857
$ws.Range("M15").Value = -8.333333333333
$ws.Range("N15").Value = -52.173913043478

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -24.242424242424
$ws.Range("M16").Value = -44.444444444444
$ws.Range("N16").Value = -85.955056179775

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -9.523809523809
$ws.Range("I17").Value = 135
$ws.Range("J17").Value = 175
$ws.Range("K17").Value = -22.857142857142
$ws.Range("L17").Value = 7.142857142857
$ws.Range("M17").Value = 27.358490566037
$ws.Range("N17").Value = -40.789473684210

# Row 18
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -70
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -51.282051282051
$ws.Range("N18").Value = -93.559322033898

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -23.076923076923
$ws.Range("I19").Value = 91
$ws.Range("J19").Value = 123
$ws.Range("K19").Value = -26.016260162601
$ws.Range("L19").Value = -35.460992907801
$ws.Range("M19").Value = -27.777777777777
$ws.Range("N19").Value = -47.093023255814

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 34
$ws.Range("L20").Value = -61.363636363636
$ws.Range("M20").Value = -71.666666666666
$ws.Range("N20").Value = -96.421052631578

# Row 21
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -11.538461538461
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = -19
$ws.Range("I21").Value = 344
$ws.Range("J21").Value = 429
$ws.Range("K21").Value = -19.813519813519
$ws.Range("L21").Value = -23.042505592841
$ws.Range("M21").Value = -27.426160337552
$ws.Range("N21").Value = -81.425485961123

# Row 23
$ws.Range("J23").Value = 29
$ws.Range("K23").Value = -34.482758620689
$ws.Range("L23").Value = -40.625
$ws.Range("M23").Value = 35.714285714285

# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = -12.765957446808
$ws.Range("I24").Value = 442
$ws.Range("J24").Value = 477
$ws.Range("K24").Value = -7.337526205450
$ws.Range("L24").Value = 7.021791767554
$ws.Range("M24").Value = -1.118568232662

# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 14.634146341463
$ws.Range("I25").Value = 245
$ws.Range("J25").Value = 225
$ws.Range("K25").Value = 8.888888888888
$ws.Range("L25").Value = 54.088050314465

# Row 26
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -19.047619047619
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 74
$ws.Range("H26").Value = -31.081081081081
$ws.Range("I26").Value = 250
$ws.Range("J26").Value = 235
$ws.Range("K26").Value = 6.382978723404
$ws.Range("L26").Value = 7.296137339055
$ws.Range("M26").Value = -36.061381074168

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 5.882352941176
$ws.Range("L27").Value = 100

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 10
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = 6.896551724137
$ws.Range("L28").Value = 40.909090909090

# Row 29
$ws.Range("M29").Value = -84.615384615384

# Row 30
$ws.Range("M30").Value = -81.818181818181

# Row 31
$ws.Range("F31").Value = 2
$ws.Range("I31").Value = 2
$ws.Range("K31").Value = 0

# ------------------------------------------------------------------
# Cells that flip between a numeric value and the "N/A" placeholder
# text (shared strings "0" / "***.*") need both their value and the
# underlying cell style (number format) updated. Value is set first
# (coercing the correct literal/text type), then the number format
# is copied from a same-column neighbour already in the target state.
# ------------------------------------------------------------------
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("D23").Value = 3
$ws.Range("D24").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = -100
$ws.Range("E24").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("D28").Value = 1
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 0
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("C31").Value = 1
$ws.Range("C24").Copy()
$ws.Range("C31").PasteSpecial(-4122)

